$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing date entries (rows 2-4)
$ws.Range("A2").Value = 45325
$ws.Range("A3").Value = 45327
$ws.Range("A4").Value = 45335

# Add new row 5: date (copy formatting from A4) + hours
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 45343
$ws.Range("B5").Value = 1.5

# Add Total column: header + SUM formula
$ws.Range("D1").Value = "Total:"
$ws.Range("D2").Formula = "=SUM(B2:B30)"

# Update selection to match saved view state
$ws.Range("F6").Select()
